$d = $word.ActiveDocument
$bm = $d.Bookmarks("_Toc255138352")
$rng = $bm.Range
$p = $rng.Paragraphs(1)
$pr = $p.Range
$ins = $d.Range($pr.End, $pr.End)
$xmlStr = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Rapid mix serves the function of distributing alum evenly throughout the raw water. The goal is to achieve small-scale mixing on the molecular scale to ensure that the alum permeates throughout all of the water coming in to the system to allow for optimal flocculation. The rapid mix system is a pipe with two orifice interfaces – one interface for macro-scale mixing, which will have </w:t></w:r><w:r><w:t>larger diameter orifices and less energy dissipation, and an interface for micro-scale mixing, which will have smaller diameter orifices and more energy dissipation, as well possibly having multiple orifices on the interface instead of just one hole. Orifices are circular, because of their ease in creation.</w:t></w:r><w:r><w:t xml:space="preserve"> The rapid mix programs serves to deliver to the user the required pipe size, the diameter of the macro-mixing and micro-mixing orifices, and the number of micro-mixing orifices.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>The diameter of the pipe used for rapid mix is determined based on the variables Q.plant, the plant flow rate, PipeSpec, which is a pipe series, and constraints on the total head loss desired, h.total, and the macromixing head loss, MacroMHconstraint. For this function, assumptions are made to the length of the pipe, L.pipe, the minor loss coefficient, K.total, and NU, the kinematic viscosity.</w:t></w:r><w:r><w:t xml:space="preserve"> In the current coding, head loss for the macro-mixing is set to 5 cm (or it can be set to a lesser value), and the total head loss is calculated based on the proposal that head loss throughout the entire plant be 40 cm.</w:t></w:r><w:r><w:t xml:space="preserve"> As such, the h.total must be determined from a total loss of 40 cm in the plant, less head loss that occurs in other parts of the plant. The head loss for the micro-mixing orifice will thus be this value, less the 5 cm of loss that is set for the macro-mixing orifice.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>h</w:t></w:r><w:r><w:t>.total for Rapid Mixer System</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>= 40cm - (</w:t></w:r><w:r><w:t>HL.Floc</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>+</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>HL.SedLaunder</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>+</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>HL.SedWeirInlet</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>+</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>HL.SedWeirExit</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>= 8.15 cm</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">When sizing the micro-mixing orifice, an energy dissipation of 1 W/kg must be met to ensure sufficient mixing.  </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>(INCLUDE MORE ON THE FUNCTIONS USED TO DETERMINE DIAMETER OF ORIFICES)</w:t></w:r></w:p>'
$ins.InsertXML($xmlStr)
Write-Output "done"
Write-Output $d.Paragraphs.Count
